$wb = $excel.ActiveWorkbook

# Sheet ALC, row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 490
$ws.Range("I18").Value = 490
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 490
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -206

# Sheet ALC, row 38
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1069.3636
$ws.Range("I38").Value = 82.875
$ws.Range("J38").Value = 3700
$ws.Range("K38").Value = 248.625
$ws.Range("L38").Value = 11100
$ws.Range("M38").Value = 123.375
$ws.Range("N38").Value = -11844

# Sheet ALC, row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1427
$ws.Range("I43").Value = 3333
$ws.Range("J43").Value = 1154.7142
$ws.Range("K43").Value = 3333
$ws.Range("L43").Value = 1154.7142
$ws.Range("M43").Value = -3264
$ws.Range("N43").Value = -1292.7142

# Sheet ALC, row 97
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 1209.75
$ws.Range("I97").Value = 3333
$ws.Range("J97").Value = 906.4286
$ws.Range("K97").Value = 9999
$ws.Range("L97").Value = 2719.2858
$ws.Range("M97").Value = -9503
$ws.Range("N97").Value = -3711.2858

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1594.5862
$ws.Range("I137").Value = 1135
$ws.Range("J137").Value = 3800.6
$ws.Range("K137").Value = 3405
$ws.Range("L137").Value = 11401.8
$ws.Range("M137").Value = -855
$ws.Range("N137").Value = -16501.8

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1695.85
$ws.Range("I61").Value = 1466.8823
$ws.Range("J61").Value = 2993.3333
$ws.Range("K61").Value = 1466.8823
$ws.Range("L61").Value = 2993.3333
$ws.Range("M61").Value = -1254.8823
$ws.Range("N61").Value = -3417.3333

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 13237892
$ws.Range("I74").Value = 19567682
$ws.Range("J74").Value = 2876.9092
$ws.Range("K74").Value = 19567682
$ws.Range("L74").Value = 2876.9092
$ws.Range("M74").Value = -19566808
$ws.Range("N74").Value = -4624.9092

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 13237892
$ws.Range("I77").Value = 19567682
$ws.Range("J77").Value = 2876.9092
$ws.Range("K77").Value = 97838410
$ws.Range("L77").Value = 14384.546
$ws.Range("M77").Value = -97834042
$ws.Range("N77").Value = -23120.546

# Sheet ARM, row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 100002350
$ws.Range("I88").Value = 1748
$ws.Range("J88").Value = 166669420
$ws.Range("K88").Value = 1748
$ws.Range("L88").Value = 166669420
$ws.Range("M88").Value = -1342
$ws.Range("N88").Value = -166670232

# Sheet ARM, row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 100002350
$ws.Range("I91").Value = 1748
$ws.Range("J91").Value = 166669420
$ws.Range("K91").Value = 1748
$ws.Range("L91").Value = 166669420
$ws.Range("M91").Value = -344
$ws.Range("N91").Value = -166672228

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3150.2727
$ws.Range("I132").Value = 3112.8076
$ws.Range("J132").Value = 3289.4285
$ws.Range("K132").Value = 9338.4228
$ws.Range("L132").Value = 9868.2855
$ws.Range("M132").Value = -6808.4228
$ws.Range("N132").Value = -14928.2855

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1695.85
$ws.Range("I136").Value = 1466.8823
$ws.Range("J136").Value = 2993.3333
$ws.Range("K136").Value = 4400.6469
$ws.Range("L136").Value = 8979.999899999999
$ws.Range("M136").Value = -1850.6469
$ws.Range("N136").Value = -14079.9999

# Sheet BSM, row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 22224838
$ws.Range("I86").Value = 40002240
$ws.Range("J86").Value = 3087.75
$ws.Range("K86").Value = 40002240
$ws.Range("L86").Value = 3087.75
$ws.Range("M86").Value = -40001117
$ws.Range("N86").Value = -5333.75

# Sheet BSM, row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 22224838
$ws.Range("I89").Value = 40002240
$ws.Range("J89").Value = 3087.75
$ws.Range("K89").Value = 200011200
$ws.Range("L89").Value = 15438.75
$ws.Range("M89").Value = -200005584
$ws.Range("N89").Value = -26670.75

# Sheet BSM, row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1426
$ws.Range("I99").Value = 1308.1818
$ws.Range("J99").Value = 1750
$ws.Range("K99").Value = 1308.1818
$ws.Range("L99").Value = 1750
$ws.Range("M99").Value = 189.8181999999999
$ws.Range("N99").Value = -4746

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3361.3125
$ws.Range("I134").Value = 2307.4443
$ws.Range("J134").Value = 4716.2856
$ws.Range("K134").Value = 6922.3329
$ws.Range("L134").Value = 14148.8568
$ws.Range("M134").Value = -4387.3329
$ws.Range("N134").Value = -19218.8568

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6170456.5
$ws.Range("I31").Value = 4370476.5
$ws.Range("J31").Value = 10528302
$ws.Range("K31").Value = 4370476.5
$ws.Range("L31").Value = 10528302
$ws.Range("M31").Value = -4370181.5
$ws.Range("N31").Value = -10528892

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6170456.5
$ws.Range("I34").Value = 4370476.5
$ws.Range("J34").Value = 10528302
$ws.Range("K34").Value = 4370476.5
$ws.Range("L34").Value = 10528302
$ws.Range("M34").Value = -4370274.5
$ws.Range("N34").Value = -10528706

# Sheet CRP, row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1224.9412
$ws.Range("I58").Value = 780.5833
$ws.Range("J58").Value = 2291.4
$ws.Range("K58").Value = 780.5833
$ws.Range("L58").Value = 2291.4
$ws.Range("M58").Value = -577.5833
$ws.Range("N58").Value = -2697.4

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2830
$ws.Range("I132").Value = 1126
$ws.Range("J132").Value = 3682
$ws.Range("K132").Value = 3378
$ws.Range("L132").Value = 11046
$ws.Range("M132").Value = -848
$ws.Range("N132").Value = -16106

# Sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5967
$ws.Range("I134").Value = 6673.1113
$ws.Range("J134").Value = 4378.25
$ws.Range("K134").Value = 20019.3339
$ws.Range("L134").Value = 13134.75
$ws.Range("M134").Value = -17484.3339
$ws.Range("N134").Value = -18204.75

# Sheet CRP, row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1224.9412
$ws.Range("I136").Value = 780.5833
$ws.Range("J136").Value = 2291.4
$ws.Range("K136").Value = 2341.7499
$ws.Range("L136").Value = 6874.200000000001
$ws.Range("M136").Value = 208.2501000000002
$ws.Range("N136").Value = -11974.2

# Sheet CUL, row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3000
$ws.Range("I5").Value = 3000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 9000
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -8888

# Sheet CUL, row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 3000
$ws.Range("I135").Value = 3000
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 27000
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -24465

# Sheet GSM, row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1876
$ws.Range("I102").Value = 1938.1177
$ws.Range("J102").Value = 1700
$ws.Range("K102").Value = 1938.1177
$ws.Range("L102").Value = 1700
$ws.Range("M102").Value = -316.1177
$ws.Range("N102").Value = -4944

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2176
$ws.Range("I132").Value = 1644.5
$ws.Range("J132").Value = 3416.1667
$ws.Range("K132").Value = 4933.5
$ws.Range("L132").Value = 10248.5001
$ws.Range("M132").Value = -2403.5
$ws.Range("N132").Value = -15308.5001

# Sheet LTW, row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 239.21739
$ws.Range("I22").Value = 224.2
$ws.Range("J22").Value = 250.76923
$ws.Range("K22").Value = 224.2
$ws.Range("L22").Value = 250.76923
$ws.Range("M22").Value = 70.80000000000001
$ws.Range("N22").Value = -840.76923

# Sheet LTW, row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 239.21739
$ws.Range("I27").Value = 224.2
$ws.Range("J27").Value = 250.76923
$ws.Range("K27").Value = 224.2
$ws.Range("L27").Value = 250.76923
$ws.Range("M27").Value = -117.2
$ws.Range("N27").Value = -464.76923

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4777
$ws.Range("I122").Value = 3605.6667
$ws.Range("J122").Value = 6785
$ws.Range("K122").Value = 10817.0001
$ws.Range("L122").Value = 20355
$ws.Range("M122").Value = -8367.000100000001
$ws.Range("N122").Value = -25255

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 20844130
$ws.Range("I132").Value = 29618500
$ws.Range("J132").Value = 5000.375
$ws.Range("K132").Value = 88855500
$ws.Range("L132").Value = 15001.125
$ws.Range("M132").Value = -88852970
$ws.Range("N132").Value = -20061.125

# Sheet LTW, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 15641464
$ws.Range("I136").Value = 25024640
$ws.Range("J136").Value = 2835
$ws.Range("K136").Value = 75073920
$ws.Range("L136").Value = 8505
$ws.Range("M136").Value = -75071370
$ws.Range("N136").Value = -13605

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 76924890
$ws.Range("I122").Value = 100001700
$ws.Range("J122").Value = 2200
$ws.Range("K122").Value = 300005100
$ws.Range("L122").Value = 6600
$ws.Range("M122").Value = -300002650
$ws.Range("N122").Value = -11500

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3349.6191
$ws.Range("I132").Value = 2714.0588
$ws.Range("J132").Value = 6050.75
$ws.Range("K132").Value = 8142.176399999999
$ws.Range("L132").Value = 18152.25
$ws.Range("M132").Value = -5612.176399999999
$ws.Range("N132").Value = -23212.25

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1661.6364
$ws.Range("I136").Value = 1264.7778
$ws.Range("J136").Value = 3447.5
$ws.Range("K136").Value = 3794.3334
$ws.Range("L136").Value = 10342.5
$ws.Range("M136").Value = -1244.3334
$ws.Range("N136").Value = -15442.5
